# Insert a new data row before row 649, shifting all following rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(649).Insert()

# Populate the newly inserted row 649 with the new record's data.
$ws.Range("A649").Value2 = 5
$ws.Range("B649").Value2 = "Macroferia Regional de Talca"
$ws.Range("C649").Value2 = "Maule"
$ws.Range("D649").Value2 = 44951
$ws.Range("E649").Value2 = 7
$ws.Range("F649").Value2 = "Fruta"
$ws.Range("G649").Value2 = 100101
$ws.Range("H649").Value2 = "Berries"
$ws.Range("I649").Value2 = 100112025
$ws.Range("J649").Value2 = "Frutilla"
$ws.Range("K649").Value2 = "Sin especificar"
$ws.Range("L649").Value2 = "Primera"
$ws.Range("M649").Value2 = 100
$ws.Range("N649").Value2 = 8000
$ws.Range("O649").Value2 = 8000
$ws.Range("P649").Value2 = 8000
$ws.Range("Q649").Value2 = "`$/caja 7 kilos"
$ws.Range("R649").Value2 = "Región del Maule"
$ws.Range("S649").Value2 = 1143
$ws.Range("T649").Value2 = 7
